$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 7).Value = 0.955666
$ws.Cells.Item(2, 8).Value = 2.866998
$ws.Cells.Item(2, 9).Value = 0.1700560254029595
$ws.Cells.Item(2, 10).Value = 0.1700560254029595
$ws.Cells.Item(2, 13).Value = 1.116695
$ws.Cells.Item(2, 14).Value = 3.350085
$ws.Cells.Item(2, 15).Value = 0.008174214292497491
$ws.Cells.Item(2, 16).Value = 0.008174214292497492
$ws.Cells.Item(2, 17).Value = 1.06718744387
$ws.Cells.Item(2, 18).Value = 9.604686994830001
$ws.Cells.Item(2, 19).Value = 0.001390074393374188
$ws.Cells.Item(2, 20).Value = 0.001390074393374188

# Row 3
$ws.Cells.Item(3, 7).Value = 0.955666
$ws.Cells.Item(3, 8).Value = 2.866998
$ws.Cells.Item(3, 9).Value = 0.1700560254029595
$ws.Cells.Item(3, 10).Value = 0.1700560254029595
$ws.Cells.Item(3, 15).Value = 0.8193429796700005
$ws.Cells.Item(3, 16).Value = 0.8193429796700005
$ws.Cells.Item(3, 17).Value = 106.9696130831067
$ws.Cells.Item(3, 18).Value = 962.7265177479601
$ws.Cells.Item(3, 19).Value = 0.1393342105644982
$ws.Cells.Item(3, 20).Value = 0.1393342105644981

# Row 4
$ws.Cells.Item(4, 7).Value = 0.955666
$ws.Cells.Item(4, 8).Value = 2.866998
$ws.Cells.Item(4, 9).Value = 0.1700560254029595
$ws.Cells.Item(4, 10).Value = 0.1700560254029595
$ws.Cells.Item(4, 15).Value = 0.172482806037502
$ws.Cells.Item(4, 16).Value = 0.1724828060375021
$ws.Cells.Item(4, 17).Value = 22.51855386953
$ws.Cells.Item(4, 18).Value = 202.66698482577
$ws.Cells.Item(4, 19).Value = 0.02933174044508719
$ws.Cells.Item(4, 20).Value = 0.02933174044508719

# Row 5
$ws.Cells.Item(5, 9).Value = 0.6638424218367511
$ws.Cells.Item(5, 10).Value = 0.663842421836751
$ws.Cells.Item(5, 13).Value = 1.116695
$ws.Cells.Item(5, 14).Value = 3.350085
$ws.Cells.Item(5, 15).Value = 0.008174214292497491
$ws.Cells.Item(5, 16).Value = 0.008174214292497492
$ws.Cells.Item(5, 17).Value = 4.165946461548334
$ws.Cells.Item(5, 18).Value = 37.493518153935
$ws.Cells.Item(5, 19).Value = 0.005426390212544119
$ws.Cells.Item(5, 20).Value = 0.005426390212544119

# Row 6
$ws.Cells.Item(6, 9).Value = 0.6638424218367511
$ws.Cells.Item(6, 10).Value = 0.663842421836751
$ws.Cells.Item(6, 15).Value = 0.8193429796700005
$ws.Cells.Item(6, 16).Value = 0.8193429796700005
$ws.Cells.Item(6, 19).Value = 0.5439146279390731
$ws.Cells.Item(6, 20).Value = 0.543914627939073

# Row 7
$ws.Cells.Item(7, 9).Value = 0.6638424218367511
$ws.Cells.Item(7, 10).Value = 0.663842421836751
$ws.Cells.Item(7, 15).Value = 0.172482806037502
$ws.Cells.Item(7, 16).Value = 0.1724828060375021
$ws.Cells.Item(7, 19).Value = 0.114501403685134
$ws.Cells.Item(7, 20).Value = 0.114501403685134

# Row 8
$ws.Cells.Item(8, 7).Value = 0.9334429999999999
$ws.Cells.Item(8, 9).Value = 0.1661015527602894
$ws.Cells.Item(8, 10).Value = 0.1661015527602894
$ws.Cells.Item(8, 13).Value = 1.116695
$ws.Cells.Item(8, 14).Value = 3.350085
$ws.Cells.Item(8, 15).Value = 0.008174214292497491
$ws.Cells.Item(8, 16).Value = 0.008174214292497492
$ws.Cells.Item(8, 17).Value = 1.042371130885
$ws.Cells.Item(8, 18).Value = 9.381340177964999
$ws.Cells.Item(8, 19).Value = 0.001357749686579184
$ws.Cells.Item(8, 20).Value = 0.001357749686579184

# Row 9
$ws.Cells.Item(9, 7).Value = 0.9334429999999999
$ws.Cells.Item(9, 9).Value = 0.1661015527602894
$ws.Cells.Item(9, 10).Value = 0.1661015527602894
$ws.Cells.Item(9, 15).Value = 0.8193429796700005
$ws.Cells.Item(9, 16).Value = 0.8193429796700005
$ws.Cells.Item(9, 18).Value = 940.3393328905798
$ws.Cells.Item(9, 19).Value = 0.1360941411664293
$ws.Cells.Item(9, 20).Value = 0.1360941411664293

# Row 10
$ws.Cells.Item(10, 7).Value = 0.9334429999999999
$ws.Cells.Item(10, 9).Value = 0.1661015527602894
$ws.Cells.Item(10, 10).Value = 0.1661015527602894
$ws.Cells.Item(10, 15).Value = 0.172482806037502
$ws.Cells.Item(10, 16).Value = 0.1724828060375021
$ws.Cells.Item(10, 19).Value = 0.02864966190728091
$ws.Cells.Item(10, 20).Value = 0.02864966190728091
